$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "foo"
$ws.Range("O1").Value = "custom column with spaces"

$ws.Range("O10").Select()
